# Update TPM-derived values in the LR-pairs sheet (Ccl5-Ackr2)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ECs -> FAPs)
$ws.Range("G2").Value = 0.3410236666666666
$ws.Range("H2").Value = 1.023071
$ws.Range("I2").Value = 0.01850325494520333
$ws.Range("J2").Value = 0.01850325494520333
$ws.Range("Q2").Value = 0.1865019783651111
$ws.Range("R2").Value = 1.678517805286
$ws.Range("S2").Value = 0.01850325494520333
$ws.Range("T2").Value = 0.01850325494520333

# Row 3 (FAPs -> FAPs)
$ws.Range("I3").Value = 0.2085050756621187
$ws.Range("J3").Value = 0.2085050756621187
$ws.Range("S3").Value = 0.2085050756621187
$ws.Range("T3").Value = 0.2085050756621187

# Row 4 (MuSCs -> FAPs)
$ws.Range("G4").Value = 0.2092423333333333
$ws.Range("H4").Value = 0.627727
$ws.Range("I4").Value = 0.0113530661283407
$ws.Range("J4").Value = 0.0113530661283407
$ws.Range("Q4").Value = 0.1144322606868889
$ws.Range("R4").Value = 1.029890346182
$ws.Range("S4").Value = 0.0113530661283407
$ws.Range("T4").Value = 0.0113530661283407

# Row 5 (Resolving-Mac -> FAPs)
$ws.Range("G5").Value = 14.03735666666667
$ws.Range("H5").Value = 42.11207
$ws.Range("I5").Value = 0.7616386032643372
$ws.Range("J5").Value = 0.7616386032643372
$ws.Range("Q5").Value = 7.676871270957778
$ws.Range("R5").Value = 69.09184143862001
$ws.Range("S5").Value = 0.7616386032643372
$ws.Range("T5").Value = 0.7616386032643372
